$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the two new columns, matching the header style used in H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill columns I (always 1) and J (mirrors column H) for data rows 2-30
for ($r = 2; $r -le 30; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
